# Elimina EC anteriores y se agregan nuevos, se modifica base de datos
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New account-statement (Estado de Cuenta) rows for the four workers,
# replacing the previous set of periods with the new ones (descending,
# most recent period 2406 first, each worker's most-recent-period row
# carries the reduced "Valor Mora" of 18560, the rest 46400).
$data = @(
    @{Row=16; Doc="1049828873"; Name="LUIS GUILLERMO CANTILLO FLOREZ";    Periodo="2406"; Valor=18560},
    @{Row=17; Doc="1049828873"; Name="LUIS GUILLERMO CANTILLO FLOREZ";    Periodo="2405"; Valor=46400},
    @{Row=18; Doc="1049828873"; Name="LUIS GUILLERMO CANTILLO FLOREZ";    Periodo="2404"; Valor=46400},
    @{Row=19; Doc="1049828873"; Name="LUIS GUILLERMO CANTILLO FLOREZ";    Periodo="2403"; Valor=46400},
    @{Row=20; Doc="1049828873"; Name="LUIS GUILLERMO CANTILLO FLOREZ";    Periodo="2402"; Valor=46400},
    @{Row=21; Doc="1049828873"; Name="LUIS GUILLERMO CANTILLO FLOREZ";    Periodo="2401"; Valor=46400},
    @{Row=22; Doc="1049828873"; Name="LUIS GUILLERMO CANTILLO FLOREZ";    Periodo="2312"; Valor=46400},
    @{Row=23; Doc="1049828873"; Name="LUIS GUILLERMO CANTILLO FLOREZ";    Periodo="2309"; Valor=46400},
    @{Row=24; Doc="1143329466"; Name="HERIBERTO DE JESUS MARQUEZ SOLIS";  Periodo="2406"; Valor=18560},
    @{Row=25; Doc="1143329466"; Name="HERIBERTO DE JESUS MARQUEZ SOLIS";  Periodo="2405"; Valor=46400},
    @{Row=26; Doc="1143329466"; Name="HERIBERTO DE JESUS MARQUEZ SOLIS";  Periodo="2404"; Valor=46400},
    @{Row=27; Doc="1143329466"; Name="HERIBERTO DE JESUS MARQUEZ SOLIS";  Periodo="2403"; Valor=46400},
    @{Row=28; Doc="1143329466"; Name="HERIBERTO DE JESUS MARQUEZ SOLIS";  Periodo="2402"; Valor=46400},
    @{Row=29; Doc="1143329466"; Name="HERIBERTO DE JESUS MARQUEZ SOLIS";  Periodo="2401"; Valor=46400},
    @{Row=30; Doc="1143329466"; Name="HERIBERTO DE JESUS MARQUEZ SOLIS";  Periodo="2312"; Valor=46400},
    @{Row=31; Doc="23139700";   Name="YESEIRA SUAREZ BATISTA";            Periodo="2406"; Valor=18560},
    @{Row=32; Doc="23139700";   Name="YESEIRA SUAREZ BATISTA";            Periodo="2405"; Valor=46400},
    @{Row=33; Doc="23139700";   Name="YESEIRA SUAREZ BATISTA";            Periodo="2404"; Valor=46400},
    @{Row=34; Doc="23139700";   Name="YESEIRA SUAREZ BATISTA";            Periodo="2403"; Valor=46400},
    @{Row=35; Doc="23139700";   Name="YESEIRA SUAREZ BATISTA";            Periodo="2402"; Valor=46400},
    @{Row=36; Doc="23139700";   Name="YESEIRA SUAREZ BATISTA";            Periodo="2401"; Valor=46400},
    @{Row=37; Doc="23139700";   Name="YESEIRA SUAREZ BATISTA";            Periodo="2312"; Valor=46400},
    @{Row=38; Doc="23139700";   Name="YESEIRA SUAREZ BATISTA";            Periodo="2311"; Valor=46400},
    @{Row=39; Doc="23139700";   Name="YESEIRA SUAREZ BATISTA";            Periodo="2308"; Valor=46400},
    @{Row=40; Doc="1051889025"; Name="MARILUZ CONEO JIMENEZ";             Periodo="2406"; Valor=18560},
    @{Row=41; Doc="1051889025"; Name="MARILUZ CONEO JIMENEZ";             Periodo="2405"; Valor=46400},
    @{Row=42; Doc="1051889025"; Name="MARILUZ CONEO JIMENEZ";             Periodo="2404"; Valor=46400},
    @{Row=43; Doc="1051889025"; Name="MARILUZ CONEO JIMENEZ";             Periodo="2403"; Valor=46400},
    @{Row=44; Doc="1051889025"; Name="MARILUZ CONEO JIMENEZ";             Periodo="2402"; Valor=46400},
    @{Row=45; Doc="1051889025"; Name="MARILUZ CONEO JIMENEZ";             Periodo="2401"; Valor=46400},
    @{Row=46; Doc="1051889025"; Name="MARILUZ CONEO JIMENEZ";             Periodo="2312"; Valor=46400},
    @{Row=47; Doc="1051889025"; Name="MARILUZ CONEO JIMENEZ";             Periodo="2311"; Valor=46400},
    @{Row=48; Doc="1051889025"; Name="MARILUZ CONEO JIMENEZ";             Periodo="2310"; Valor=46400}
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 3).Value = $item.Doc
    $ws.Cells.Item($r, 4).Value = $item.Name
    $ws.Cells.Item($r, 5).Value = $item.Periodo
    $ws.Cells.Item($r, 6).Value = $item.Valor
}
